# DSS.xlsx - "Add files via upload" edit
#
# Two new trainees are appended to the training-record sheet:
#   - Mohamed Abdel Hamid Basyouni Agiza  (certs DSS2302-DSS2309) -> rows 1303-1310
#   - Hany Mohamed Salah Eldin Ahmed      (certs DSS2343-DSS2350) -> rows 1344-1351
#     plus a First Aid record (DSS2367, 15-3-2025)                -> row 1368
#
# Rows 1303-1310 and 1344-1351 already exist in the sheet as blank placeholder
# rows (with the correct cell styles already applied) - we only need to fill
# in their values. Rows 1352-1367 need to be turned from single-cell
# placeholder rows into full 5-column blank rows matching the style of the
# row above them, and row 1368 is a brand new data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Block 1: rows 1303-1310 - Mohamed Abdel Hamid Basyouni Agiza
# ---------------------------------------------------------------------------
# Name first (single new shared string), then the certificate numbers in
# order, then the (already-existing) course names / dates / the "1" flag
# copied down from the identical course block above (rows 1295-1302).

$ws.Range("B1303").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1304").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1305").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1306").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1307").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1308").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1309").Value = "Mohamed Abdel Hamid Basyouni Agiza"
$ws.Range("B1310").Value = "Mohamed Abdel Hamid Basyouni Agiza"

$ws.Range("A1303").Value = "DSS2302"
$ws.Range("A1304").Value = "DSS2303"
$ws.Range("A1305").Value = "DSS2304"
$ws.Range("A1306").Value = "DSS2305"
$ws.Range("A1307").Value = "DSS2306"
$ws.Range("A1308").Value = "DSS2307"
$ws.Range("A1309").Value = "DSS2308"
$ws.Range("A1310").Value = "DSS2309"

$ws.Range("C1303").Value = $ws.Range("C1295").Value2
$ws.Range("C1304").Value = $ws.Range("C1296").Value2
$ws.Range("C1305").Value = $ws.Range("C1297").Value2
$ws.Range("C1306").Value = $ws.Range("C1298").Value2
$ws.Range("C1307").Value = $ws.Range("C1299").Value2
$ws.Range("C1308").Value = $ws.Range("C1300").Value2
$ws.Range("C1309").Value = $ws.Range("C1301").Value2
$ws.Range("C1310").Value = $ws.Range("C1302").Value2

$ws.Range("D1303").Value = 45779
$ws.Range("D1304").Value = 45932
$ws.Range("D1305").Value = 45810
$ws.Range("D1306").Value = 45718
$ws.Range("D1307").Value = 45659
$ws.Range("D1308").Value = 45690
$ws.Range("D1309").Value = 45840
$ws.Range("D1310").Value = 45871

$ws.Range("E1303").Value = 1
$ws.Range("E1304").Value = 1
$ws.Range("E1305").Value = 1
$ws.Range("E1306").Value = 1
$ws.Range("E1307").Value = 1
$ws.Range("E1308").Value = 1
$ws.Range("E1309").Value = 1
$ws.Range("E1310").Value = 1

# ---------------------------------------------------------------------------
# Block 2: rows 1344-1351 - Hany Mohamed Salah Eldin Ahmed
# ---------------------------------------------------------------------------
# These rows are currently single-cell (column A only) placeholder rows, so
# we first stamp the A1344:E1351 block with the formatting used by the
# identical course block above (rows 1295-1302) and only then write values -
# writing the values first means PasteSpecial would leave them untouched,
# writing them after the format copy keeps the original style indices
# (s="44" / s="17") instead of minting new ones.

$ws.Range("A1295:E1302").Copy() | Out-Null
$ws.Range("A1344:E1351").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B1344").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1345").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1346").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1347").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1348").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1349").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1350").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("B1351").Value = "Hany Mohamed Salah Eldin Ahmed"

$ws.Range("A1344").Value = "DSS2343"
$ws.Range("A1345").Value = "DSS2344"
$ws.Range("A1346").Value = "DSS2345"
$ws.Range("A1347").Value = "DSS2346"
$ws.Range("A1348").Value = "DSS2347"
$ws.Range("A1349").Value = "DSS2348"
$ws.Range("A1350").Value = "DSS2349"
$ws.Range("A1351").Value = "DSS2350"

$ws.Range("C1344").Value = $ws.Range("C1295").Value2
$ws.Range("C1345").Value = $ws.Range("C1296").Value2
$ws.Range("C1346").Value = $ws.Range("C1297").Value2
$ws.Range("C1347").Value = $ws.Range("C1298").Value2
$ws.Range("C1348").Value = $ws.Range("C1299").Value2
$ws.Range("C1349").Value = $ws.Range("C1300").Value2
$ws.Range("C1350").Value = $ws.Range("C1301").Value2
$ws.Range("C1351").Value = $ws.Range("C1302").Value2

$ws.Range("D1344").Value = 45779
$ws.Range("D1345").Value = 45932
$ws.Range("D1346").Value = 45810
$ws.Range("D1347").Value = 45718
$ws.Range("D1348").Value = 45659
$ws.Range("D1349").Value = 45690
$ws.Range("D1350").Value = 45840
$ws.Range("D1351").Value = 45871

$ws.Range("E1344").Value = 1
$ws.Range("E1345").Value = 1
$ws.Range("E1346").Value = 1
$ws.Range("E1347").Value = 1
$ws.Range("E1348").Value = 1
$ws.Range("E1349").Value = 1
$ws.Range("E1350").Value = 1
$ws.Range("E1351").Value = 1

# ---------------------------------------------------------------------------
# Rows 1352-1367: blank spacer rows, re-formatted (A column only, s="44")
# ---------------------------------------------------------------------------

$ws.Range("A1344").Copy() | Out-Null
$ws.Range("A1352:A1367").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 1368: Hany Mohamed Salah Eldin Ahmed - First Aid - 15-3-2025
# ---------------------------------------------------------------------------

$ws.Range("A1342:E1342").Copy() | Out-Null
$ws.Range("A1368:E1368").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D1368").Value = "15-3-2025"
$ws.Range("A1368").Value = "DSS2367"
$ws.Range("B1368").Value = "Hany Mohamed Salah Eldin Ahmed"
$ws.Range("C1368").Value = $ws.Range("C1342").Value2
$ws.Range("E1368").Value = 1

# ---------------------------------------------------------------------------
# Manual page break: was after row 1344, now after row 1343 (the new trainee
# block now starts cleanly at the top of its own printed page).
# ---------------------------------------------------------------------------

$ws.HPageBreaks.Item($ws.HPageBreaks.Count).Delete()
$ws.HPageBreaks.Add($ws.Range("A1344")) | Out-Null

# ---------------------------------------------------------------------------
# Scroll position / selection, matching the author's last on-screen state.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 1296
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J1300").Select() | Out-Null
